# Update the roster table: swap a few players (and their position/team)
# for new ones. Rows 3,4,5,7,8,9,12,13,16,17 keep their original
# player/position/team; rows 2,6,10,11,14,15 are replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Tyrese Maxey"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Philadelphia 76ers"

$ws.Range("A6").Value = "Ausar Thompson"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Detroit Pistons"

$ws.Range("A10").Value = "Jalen Duren"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Detroit Pistons"

$ws.Range("A11").Value = "Jarrett Allen"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Cleveland Cavaliers"

$ws.Range("A14").Value = "Franz Wagner"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "Orlando Magic"

$ws.Range("A15").Value = "Darius Garland"
$ws.Range("B15").Value = "PG"
$ws.Range("C15").Value = "Cleveland Cavaliers"
